# update scripts wuth new tpm
# Re-run of the NATMI TPM pipeline for Wnt5a-Fzd4 (YoungD2) updated the
# per-pair TPM-derived statistics and dropped the "Resolving-Mac" target
# cluster, shrinking the LR-pair table from 10 data rows to 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-5: refresh metrics recomputed from the new TPM input (ligand/receptor specificity, edge weights) ---
# Row 2
$ws.Cells.Item(2,9).Value = 0.977669497583861
$ws.Cells.Item(2,10).Value = 0.977669497583861
$ws.Cells.Item(2,13).Value = 26.60444266666667
$ws.Cells.Item(2,14).Value = 79.813328
$ws.Cells.Item(2,15).Value = 0.5736225649467147
$ws.Cells.Item(2,16).Value = 0.5736225649467147
$ws.Cells.Item(2,17).Value = 199.9269327292533
$ws.Cells.Item(2,18).Value = 1799.34239456328
$ws.Cells.Item(2,19).Value = 0.5608132848742202
$ws.Cells.Item(2,20).Value = 0.5608132848742202

# Row 3
$ws.Cells.Item(3,9).Value = 0.977669497583861
$ws.Cells.Item(3,10).Value = 0.977669497583861
$ws.Cells.Item(3,15).Value = 0.3025349071358453
$ws.Cells.Item(3,16).Value = 0.3025349071358453
$ws.Cells.Item(3,18).Value = 948.993148997475
$ws.Cells.Item(3,19).Value = 0.2957791506610819
$ws.Cells.Item(3,20).Value = 0.2957791506610819

# Row 4
$ws.Cells.Item(4,9).Value = 0.977669497583861
$ws.Cells.Item(4,10).Value = 0.977669497583861
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.2022123333333333
$ws.Cells.Item(4,14).Value = 0.606637
$ws.Cells.Item(4,15).Value = 0.004359931864156574
$ws.Cells.Item(4,16).Value = 0.004359931864156574
$ws.Cells.Item(4,17).Value = 1.519584231471667
$ws.Cells.Item(4,18).Value = 13.676258083245
$ws.Cells.Item(4,19).Value = 0.004262572395129825
$ws.Cells.Item(4,20).Value = 0.004262572395129825

# Row 5
$ws.Cells.Item(5,9).Value = 0.977669497583861
$ws.Cells.Item(5,10).Value = 0.977669497583861
$ws.Cells.Item(5,13).Value = 5.541567000000001
$ws.Cells.Item(5,14).Value = 16.624701
$ws.Cells.Item(5,15).Value = 0.1194825960532834
$ws.Cells.Item(5,16).Value = 0.1194825960532834
$ws.Cells.Item(5,17).Value = 41.643739983765
$ws.Cells.Item(5,18).Value = 374.793659853885
$ws.Cells.Item(5,19).Value = 0.116814489653429
$ws.Cells.Item(5,20).Value = 0.116814489653429

# --- Rows 6-9: re-derived for the MuSCs sending cluster (target-cluster set no longer includes Resolving-Mac) ---
# Row 6
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Wnt5a"
$ws.Cells.Item(6,3).Value = "Fzd4"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.171642
$ws.Cells.Item(6,8).Value = 0.514926
$ws.Cells.Item(6,9).Value = 0.02233050241613897
$ws.Cells.Item(6,10).Value = 0.02233050241613898
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 26.60444266666667
$ws.Cells.Item(6,14).Value = 79.813328
$ws.Cells.Item(6,15).Value = 0.5736225649467147
$ws.Cells.Item(6,16).Value = 0.5736225649467147
$ws.Cells.Item(6,17).Value = 4.566439748192
$ws.Cells.Item(6,18).Value = 41.097957733728
$ws.Cells.Item(6,19).Value = 0.01280928007249445
$ws.Cells.Item(6,20).Value = 0.01280928007249445

# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Wnt5a"
$ws.Cells.Item(7,3).Value = "Fzd4"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.171642
$ws.Cells.Item(7,8).Value = 0.514926
$ws.Cells.Item(7,9).Value = 0.02233050241613897
$ws.Cells.Item(7,10).Value = 0.02233050241613898
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 14.03147833333333
$ws.Cells.Item(7,14).Value = 42.094435
$ws.Cells.Item(7,15).Value = 0.3025349071358453
$ws.Cells.Item(7,16).Value = 0.3025349071358453
$ws.Cells.Item(7,17).Value = 2.40839100409
$ws.Cells.Item(7,18).Value = 21.67551903681
$ws.Cells.Item(7,19).Value = 0.006755756474763373
$ws.Cells.Item(7,20).Value = 0.006755756474763374

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Wnt5a"
$ws.Cells.Item(8,3).Value = "Fzd4"
$ws.Cells.Item(8,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.171642
$ws.Cells.Item(8,8).Value = 0.514926
$ws.Cells.Item(8,9).Value = 0.02233050241613897
$ws.Cells.Item(8,10).Value = 0.02233050241613898
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.2022123333333333
$ws.Cells.Item(8,14).Value = 0.606637
$ws.Cells.Item(8,15).Value = 0.004359931864156574
$ws.Cells.Item(8,16).Value = 0.004359931864156574
$ws.Cells.Item(8,17).Value = 0.034708129318
$ws.Cells.Item(8,18).Value = 0.312373163862
$ws.Cells.Item(8,19).Value = 0.00009735946902674968
$ws.Cells.Item(8,20).Value = 0.00009735946902674969

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Wnt5a"
$ws.Cells.Item(9,3).Value = "Fzd4"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.171642
$ws.Cells.Item(9,8).Value = 0.514926
$ws.Cells.Item(9,9).Value = 0.02233050241613897
$ws.Cells.Item(9,10).Value = 0.02233050241613898
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 5.541567000000001
$ws.Cells.Item(9,14).Value = 16.624701
$ws.Cells.Item(9,15).Value = 0.1194825960532834
$ws.Cells.Item(9,16).Value = 0.1194825960532834
$ws.Cells.Item(9,17).Value = 0.9511656430140001
$ws.Cells.Item(9,18).Value = 8.560490787126001
$ws.Cells.Item(9,19).Value = 0.002668106399854402
$ws.Cells.Item(9,20).Value = 0.002668106399854402

# --- Rows 10-11 (old MuSCs/MuSCs and MuSCs/Resolving-Mac pairs) no longer exist; drop them ---
$ws.Rows("10:11").Delete()
